$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Clear out most of the "NT - not tested" cells in the BinSeg/SegNeigh/PELT
# columns (F:H) that were re-run and no longer apply, leaving only what was
# actually retested.
$ws.Range("H6").ClearContents()
$ws.Range("F7:H7").ClearContents()
$ws.Range("F8:H8").ClearContents()
$ws.Range("F9").ClearContents()
$ws.Range("H9").ClearContents()
$ws.Range("F10:H10").ClearContents()
$ws.Range("F11:H11").ClearContents()
$ws.Range("F13:H13").ClearContents()

# SegNeigh (column G) for the "MBIC" penalty (row 9) was tested and passed.
$ws.Range("G9").Value = "√"

# The "Asymptotic" penalty row (12) was previously marked "NA" for every
# method; it has now actually been tested and found not-tested ("NT").
$ws.Range("E12:H12").Value = "NT"

# Move the active selection.
$ws.Range("F10").Select() | Out-Null
